# Slide 17 ("Can Tetraplex be a unicorn as a SuperApp ?"): fix the
# "887.3 billion USD t by 2033 that's 1.77 billion USD (Estimated market 2033)"
# line so it reads
# "887.3 billion USD  by 2033 that's 1.77 billion USD (Estimated market by 2033)"
# and split the trailing parenthetical into its own runs, same as the
# authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Work right-to-left across the paragraph so earlier character offsets
# stay valid while later ranges are rewritten.

# "(Estimated market 2033)" (chars 88-110) -> split into four runs:
#   "(" / "Estimated " / "market by " / "2033)"
$tr.Characters(106, 5).Text  = "2033)"
$tr.Characters(99, 7).Text   = "market by "
$tr.Characters(89, 10).Text  = "Estimated "
$tr.Characters(88, 1).Text   = "("

# "t by 2033 " + "that's " (chars 54-70, two runs with identical default
# formatting) -> merge/rewrite as a single run " by 2033 that's "
$tr.Characters(54, 17).Text = " by 2033 that’s "
